$wb = $excel.ActiveWorkbook

# The F-column "想去人数" (want-to-go count) figures were refreshed on sheets
# "展览" and "全部类型" for rows 4, 7, 8, 10, 12, 15.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value = 1551
    $ws.Range("F7").Value = 11249
    $ws.Range("F8").Value = 9
    $ws.Range("F10").Value = 382
    $ws.Range("F12").Value = 1080
    $ws.Range("F15").Value = 12914
}
